$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3520279824733734
$ws.Range("B1").Value = 0.3882378935813904
$ws.Range("D1").Value = 2.421206951141357
$ws.Range("E1").Value = 1.167293071746826
